$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded workbook dropped the row-2 review (appid/captain-credit entry
# for eligitel@gmail.com / ronenchen27@gmail.com, dated 25/5/2019 13:44).
# Remove that row entirely; Excel shifts every following row up by one.
#
# The engine's native row-delete does not re-anchor the sheet's Hyperlinks
# collection (their Range stays pinned to the old row numbers), so we drop
# all existing hyperlinks first and rebuild them against the post-delete
# layout.
$ws.Hyperlinks.Delete()
$ws.Rows.Item(2).Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:zaittomer@gmail.com", "", "", "zaittomer@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:rontiddler560@gmail.com", "", "", "rontiddler560@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:gregneri12@gmail.com", "", "", "gregneri12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:snizzvered@gmail.com", "", "", "snizzvered@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:budoyoni2@gmail.com", "", "", "budoyoni2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:hermanliran@gmail.com", "", "", "hermanliran@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:gazittalia1@gmail.com", "", "", "gazittalia1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:hermanliran@gmail.com", "", "", "hermanliran@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:leviadlevi22@gmail.com", "", "", "leviadlevi22@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:gazittalia1@gmail.com", "", "", "gazittalia1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:freelancernachus@gmail.com", "", "", "freelancernachus@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:nevilgreen@gmail.com", "", "", "nevilgreen@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:stclerari834@gmail.com", "", "", "stclerari834@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:stcydouel274@gmail.com", "", "", "stcydouel274@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D17"), "mailto:sinuspai@gmail.com", "", "", "sinuspai@gmail.com")

$ws.Range("B2").Select()
